# Update industrial mapping scheme for Austria (CDM -> CDH and adapted percentages)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New multi-line mapping text for cell B2 (CR+PC and S rows changed from CDM to CDH,
# with adapted percentages)
$newText = "7% MUR/LWAL+CDN/H:1`n27% CR+PC/LFM+CDL/H:1`n1% CR+PC/LFM+CDH/H:1`n6% S+SL/LFM+CDL/H:1`n30% S/LFM+CDL/RME/H:1`n3% S/LFM+CDH/H:1`n5% W/LWAL+CDL/H:1`n19% CR/LFINF+CDL/H:2`n2% CR/LFINF+CDL/HBET:3-5"

$cell = $ws.Cells.Item(2, 2)
$cell.Value = $newText

# Wrap the (now longer) text and resize the row/column so it remains fully visible
$cell.WrapText = $true
$ws.Rows.Item(2).RowHeight = 141
$ws.Columns.Item(2).ColumnWidth = 52

# Make B2 the active / selected cell, matching the saved view state
$ws.Activate()
$ws.Range("B2").Select()
